$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_schedule_1")

# Relativize the recorded video file paths (drop the absolute
# "C:\Users\...\streamscheduler\" prefix), keeping the same
# logical file for each row.
$ws.Range("C2").Value = "test_files\vids\test.mp4"
$ws.Range("C3").Value = "test_files\vids2\test2.mp4"
$ws.Range("C4").Value = "test_files\vids\test4.mp4"

# Restore the selection that was active when the workbook was saved.
$ws.Activate()
$ws.Range("D13").Select()
